# experiments_bookkeeping.xlsx update:
#  - trim the experiment log from 16 data rows down to 5 (rows 3-7)
#  - rework which "Agent versions" flags are set (Full/Partial) and
#    which Observation/Long-horizon/Reconstruction columns (C/D/E) are marked
#  - keep the P column "ablation cost" helper formula consistent over the
#    smaller range
#  - refresh the selection / filter range / used dimension accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove the now-unused trailing rows (old rows 8-18) --------------
$ws.Rows("8:18").Delete()

# --- row 4: drop the Observation "x" mark ------------------------------
$ws.Range("D4").ClearContents()

# --- row 5: Full -> Partial, keep no Observation/Long-horizon marks ----
$ws.Range("B5").Value = "Partial"

# --- row 6: Full -> Partial; Observation mark moves out, Reconstruction
#            mark (E) is newly added, Long-horizon (D) mark kept --------
$ws.Range("B6").Value = "Partial"
$ws.Range("C6").ClearContents()
$ws.Range("E6").Value = "x"

# --- row 7: Full -> Partial; Observation mark kept, Long-horizon and
#            Reconstruction marks newly added --------------------------
$ws.Range("B7").Value = "Partial"
$ws.Range("D7").Value = "x"
$ws.Range("E7").Value = "x"

# --- re-establish the shared formula block for the smaller P4:P7 range -
$ws.Range("P4:P7").Formula = "=IF(ISBLANK(C4),0,0.01)"

# --- update the autofilter / defined name range to match the new data --
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Experiments!`$B`$1:`$U`$7"

# --- move the active selection to E7, matching the saved file ----------
$ws.Range("E7").Select()
